# "Add files via upload" — adds a new "ApPredict version information"
# worksheet (with ApPredict/Chaste build-provenance metadata) to the
# workbook, as the last/active sheet.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet after the last existing sheet ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ApPredict version information"

# --- Row data: Column A = label, Column B = value -----------------------
$rows = @(
    @("ApPredict Version", "37cc5a6"),
    @("Chaste Version", "2019.1.682dce0"),
    @("Modified", $true),
    @("Build options", "GccOpt, shared libraries"),
    @("OS info", "Linux d09b088bdc9f 4.15.0-161-generic #169-Ubuntu SMP Fri Oct 15 13:41:54 UTC 2021 x86_64"),
    @("Compiler", "gcc, version b'9.3.0'"),
    @("Compiler flags", "-O3 -std=c++14"),
    @("XSD", "4.0.0"),
    @("VTK", "no"),
    @("Xerces", "3.2.0"),
    @("SUNDIALS", "2.5.0"),
    @("HDF5", "1.8.16"),
    @("Boost", "1.65.1"),
    @("PETSc", "3.12.4"),
    @("Parmetis", "4.0.3"),
    @("Ap Predict arguments", " --pacing-freq 1 --pacing-max-time 5 --plasma-conc-high 100 --plasma-conc-low 0 --plasma-conc-count 4 --plasma-conc-logscale true --model 1")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $newSheet.Cells.Item($r, 1).Value = $rows[$i][0]

    $v = $rows[$i][1]
    if ($v -is [string] -and $v -eq "1.8.16") {
        # Plain assignment gets misread as a date (e.g. "1.8.16" -> Jan
        # 2016) by the auto-detecting Value setter. Route it through a
        # text formula + paste-special-values so it lands as a genuine
        # shared-string cell, same as every other text value here.
        $cell = $newSheet.Cells.Item($r, 2)
        $cell.Formula = '="' + $v + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    } else {
        $newSheet.Cells.Item($r, 2).Value = $v
    }
}

# --- Match the saved selection / active-sheet state ---------------------
$newSheet.Range("A1:B16").Select()
$newSheet.Activate()
